$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-10 Wednesday" "2025-12-11 Thursday"

Replace-Text "15×15=" "74×78="
Replace-Text "90×43=" "77×22="
Replace-Text "24×99=" "42×46="
Replace-Text "23×43=" "27×96="
Replace-Text "50×94=" "20×50="
Replace-Text "24×23=" "24×48="
Replace-Text "77×59=" "72×82="
Replace-Text "97×79=" "44×53="
Replace-Text "86×33=" "41×35="
Replace-Text "33×39=" "99×63="
Replace-Text "48×41=" "31×71="
Replace-Text "13×53=" "26×16="
Replace-Text "96×53=" "24×19="
Replace-Text "14×77=" "31×38="
Replace-Text "13×26=" "15×74="
Replace-Text "31×50=" "51×88="
Replace-Text "36×16=" "59×42="
Replace-Text "85×54=" "80×92="
Replace-Text "44×65=" "14×53="
Replace-Text "24×56=" "47×23="
Replace-Text "88×11=" "81×81="
Replace-Text "43×63=" "41×41="
Replace-Text "63×93=" "40×83="
Replace-Text "91×33=" "95×73="
Replace-Text "92×83=" "95×89="
